# Apply the edit described by the diff:
#  1. Within each 4-row "year block" (rows 2..69), swap the entire contents
#     (columns A:E) of the second row ("B" record) and the third row
#     ("C" record). The first ("A") and fourth ("D") rows of each block
#     stay where they are.
#  2. Delete columns F and G entirely (they duplicated data already present
#     in columns B and E and are removed from the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$blockStarts = 2,6,10,14,18,22,26,30,34,38,42,46,50,54,58,62,66

foreach ($start in $blockStarts) {
    $rowB = $start + 1
    $rowC = $start + 2

    $rangeB = "A" + $rowB + ":E" + $rowB
    $rangeC = "A" + $rowC + ":E" + $rowC

    $valuesB = $ws.Range($rangeB).Value2
    $valuesC = $ws.Range($rangeC).Value2

    $ws.Range($rangeB).Value2 = $valuesC
    $ws.Range($rangeC).Value2 = $valuesB
}

# Remove the now-obsolete columns F and G (electronics product-sales-ratio /
# sales-volume duplicate columns).
$ws.Range("F1:G69").Delete()
